$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 21 with the new "user story" item, reusing the same
# formatting (fill/style) as row 20 immediately above it.
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Ny sida när testet är skapat/alternativt meddelande "

# Update the active selection to the newly added cell.
$ws.Range("A21").Select()
